$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H: "Save", styled like the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Rows that should be marked as a "Save" (value 1); all others default to 0
$saveRows = @(31, 38)

for ($r = 2; $r -le 49; $r++) {
    if ($saveRows -contains $r) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
